$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sostdc1 -> Lrp5, Target cluster: ECs) updated TPM-derived values
$ws.Range("G2").Value = 0.6311703333333334
$ws.Range("H2").Value = 1.893511
$ws.Range("M2").Value = 15.24491733333333
$ws.Range("N2").Value = 45.73475199999999
$ws.Range("O2").Value = 0.4831257321597052
$ws.Range("P2").Value = 0.4831257321597052
$ws.Range("Q2").Value = 9.62213955491911
$ws.Range("R2").Value = 86.599255994272
$ws.Range("S2").Value = 0.4831257321597052
$ws.Range("T2").Value = 0.4831257321597052

# Row 3 (Sostdc1 -> Lrp5, Target cluster: FAPs) updated TPM-derived values
$ws.Range("G3").Value = 0.6311703333333334
$ws.Range("H3").Value = 1.893511
$ws.Range("O3").Value = 0.327710667227878
$ws.Range("P3").Value = 0.327710667227878
$ws.Range("Q3").Value = 6.526826380384001
$ws.Range("R3").Value = 58.74143742345601
$ws.Range("S3").Value = 0.327710667227878
$ws.Range("T3").Value = 0.327710667227878

# Row 4 (Sostdc1 -> Lrp5, Target cluster: MuSCs) updated TPM-derived values
$ws.Range("G4").Value = 0.6311703333333334
$ws.Range("H4").Value = 1.893511
$ws.Range("M4").Value = 5.969012333333333
$ws.Range("N4").Value = 17.907037
$ws.Range("O4").Value = 0.1891636006124168
$ws.Range("P4").Value = 0.1891636006124168
$ws.Range("Q4").Value = 3.767463504100778
$ws.Range("R4").Value = 33.907171536907
$ws.Range("S4").Value = 0.1891636006124168
$ws.Range("T4").Value = 0.1891636006124168
